$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A21:H21").Merge()
$ws.Range("A21").Value = "Executions"
$ws.Range("A21:H21").Style = "Check Cell"
